$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source added two new weekly price records for Kiwi (Terminal
# Hortofruticola Agro Chillan) right after the current header block,
# pushing the existing historical rows (old 308:347) down to 310:349.
$ws.Rows("308:309").Insert()

# New row 308: "Primera" quality, week of 2023-07-25 (serial 45132)
$ws.Cells.Item(308, 1).Value2  = 7
$ws.Cells.Item(308, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(308, 3).Value2  = "Ñuble"
$ws.Cells.Item(308, 4).Value2  = 45132
$ws.Cells.Item(308, 5).Value2  = 16
$ws.Cells.Item(308, 6).Value2  = "Fruta"
$ws.Cells.Item(308, 7).Value2  = 100101
$ws.Cells.Item(308, 8).Value2  = "Berries"
$ws.Cells.Item(308, 9).Value2  = 100101007
$ws.Cells.Item(308, 10).Value2 = "Kiwi"
$ws.Cells.Item(308, 11).Value2 = "Hayward"
$ws.Cells.Item(308, 12).Value2 = "Primera"
$ws.Cells.Item(308, 13).Value2 = 60
$ws.Cells.Item(308, 14).Value2 = 12000
$ws.Cells.Item(308, 15).Value2 = 12000
$ws.Cells.Item(308, 16).Value2 = 12000
$ws.Cells.Item(308, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(308, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(308, 19).Value2 = 667
$ws.Cells.Item(308, 20).Value2 = 18

# New row 309: "Segunda" quality, same week (serial 45132)
$ws.Cells.Item(309, 1).Value2  = 7
$ws.Cells.Item(309, 2).Value2  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(309, 3).Value2  = "Ñuble"
$ws.Cells.Item(309, 4).Value2  = 45132
$ws.Cells.Item(309, 5).Value2  = 16
$ws.Cells.Item(309, 6).Value2  = "Fruta"
$ws.Cells.Item(309, 7).Value2  = 100101
$ws.Cells.Item(309, 8).Value2  = "Berries"
$ws.Cells.Item(309, 9).Value2  = 100101007
$ws.Cells.Item(309, 10).Value2 = "Kiwi"
$ws.Cells.Item(309, 11).Value2 = "Hayward"
$ws.Cells.Item(309, 12).Value2 = "Segunda"
$ws.Cells.Item(309, 13).Value2 = 80
$ws.Cells.Item(309, 14).Value2 = 10000
$ws.Cells.Item(309, 15).Value2 = 10000
$ws.Cells.Item(309, 16).Value2 = 10000
$ws.Cells.Item(309, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(309, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(309, 19).Value2 = 556
$ws.Cells.Item(309, 20).Value2 = 18
